# Fixing a set for departement and position in xlsxReader
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Dipartimento (H) / Posizione (I) values -----------------------------
$ws.Range("H2").Value = 11
$ws.Range("I2").Value = 11

$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 3

$ws.Range("H4").Value = 6
$ws.Range("I4").Value = 5

$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 3

$ws.Range("I6").Value = 7

$ws.Range("H7").Value = 7
$ws.Range("I7").Value = 4

$ws.Range("H8").Value = 5
$ws.Range("I8").Value = 4

$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 4

$ws.Range("H10").Value = 8
$ws.Range("I10").Value = 3

$ws.Range("H11").Value = 3
$ws.Range("I11").Value = 9

$ws.Range("H12").Value = 10
$ws.Range("I12").Value = 5

$ws.Range("H13").Value = 3
$ws.Range("I13").Value = 3

$ws.Range("I14").Value = 7

$ws.Range("H15").Value = 2
$ws.Range("I15").Value = 3

$ws.Range("H16").Value = 8
$ws.Range("I16").Value = 4

$ws.Range("H17").Value = 3

$ws.Range("H18").Value = 5
$ws.Range("I18").Value = 3

$ws.Range("H19").Value = 2

$ws.Range("H20").Value = 3
$ws.Range("I20").Value = 5

$ws.Range("H21").Value = 9
$ws.Range("I21").Value = 4

# --- Swap the two Riccardo/Riccardino e-mail hyperlinks ------------------
# Row 20 (Riccardo Gatti) keeps the Hyperlink style but now correctly
# points at riccardo@gatti.it instead of riccardo@gattinone.it.
$ws.Range("E20").Hyperlinks.Delete()
$ws.Range("E20").Value = "riccardo@gatti.it"
$ws.Hyperlinks.Add($ws.Range("E20"), "mailto:riccardo@gatti.it")
$ws.Range("E20").Style = "Hyperlink"

# Row 21 (Riccardino Gattinone) now gets its own hyperlink pointing at
# riccardo@gattinone.it.
$ws.Range("E21").Value = "riccardo@gattinone.it"
$ws.Hyperlinks.Add($ws.Range("E21"), "mailto:riccardo@gattinone.it")
$ws.Range("E21").Style = "Hyperlink"

# --- New column widths for H/I -------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 10.83
$ws.Columns.Item(9).ColumnWidth = 10.0

# --- Selection / view state -----------------------------------------------
$ws.Range("I19").Select() | Out-Null
